$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guide Quests")

# Update the quest numbering in column A for rows 8-16 (21-29 -> 12-20)
$ws.Range("A8").Value = 12
$ws.Range("A9").Value = 13
$ws.Range("A10").Value = 14
$ws.Range("A11").Value = 15
$ws.Range("A12").Value = 16
$ws.Range("A13").Value = 17
$ws.Range("A14").Value = 18
$ws.Range("A15").Value = 19
$ws.Range("A16").Value = 20
